$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B1").Value2 = "0: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 1: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 2: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B2").Value2 = "3: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 4: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 5: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B3").Value2 = "6: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 7: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 8: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B4").Value2 = "9: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 10: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 11: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B5").Value2 = "12: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 13: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 14: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B6").Value2 = "15: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 16: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 17: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B7").Value2 = "18: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 19: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 20: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B8").Value2 = "21: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 22: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 23: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B9").Value2 = "24: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 25: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 26: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B10").Value2 = "27: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 28: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 29: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B11").Value2 = "30: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 31: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 32: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B12").Value2 = "33: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 34: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 35: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B13").Value2 = "36: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 37: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 38: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B14").Value2 = "39: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 40: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 41: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B15").Value2 = "42: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 43: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 44: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B16").Value2 = "45: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 46: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 47: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B17").Value2 = "48: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 49: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 50: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B18").Value2 = "51: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 52: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 53: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B19").Value2 = "54: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 55: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 56: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B20").Value2 = "57: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 58: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 59: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B21").Value2 = "60: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 61: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 62: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B22").Value2 = "63: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 64: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 65: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B23").Value2 = "66: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 67: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 68: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B24").Value2 = "69: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 70: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 71: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B25").Value2 = "72: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 73: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 74: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B26").Value2 = "75: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 76: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 77: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B27").Value2 = "78: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 79: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 80: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B28").Value2 = "81: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 82: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 83: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B29").Value2 = "84: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 85: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 86: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B30").Value2 = "87: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 88: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 89: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B31").Value2 = "90: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 91: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 92: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B32").Value2 = "93: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 94: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 95: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B33").Value2 = "96: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 97: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 98: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B34").Value2 = "99: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 100: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 101: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B35").Value2 = "102: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 103: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 104: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B36").Value2 = "105: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 106: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 107: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B37").Value2 = "108: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 109: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 110: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B38").Value2 = "111: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 112: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 113: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B39").Value2 = "114: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 115: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 116: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B40").Value2 = "117: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 118: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 119: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B41").Value2 = "120: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 121: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 122: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B42").Value2 = "123: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 124: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 125: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B43").Value2 = "126: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 127: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 128: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B44").Value2 = "129: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 130: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 131: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B45").Value2 = "132: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 133: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 134: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B46").Value2 = "135: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 136: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 137: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B47").Value2 = "138: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 139: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 140: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B48").Value2 = "141: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 142: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 143: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B49").Value2 = "144: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 145: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 146: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B50").Value2 = "147: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 148: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 149: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B51").Value2 = "150: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 151: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 152: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B52").Value2 = "153: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 154: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 155: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B53").Value2 = "156: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 157: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 158: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B54").Value2 = "159: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 160: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 161: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B55").Value2 = "162: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 163: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 164: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B56").Value2 = "165: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 166: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 167: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B57").Value2 = "168: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 169: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 170: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B58").Value2 = "171: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 172: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 173: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B59").Value2 = "174: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 175: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 176: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B60").Value2 = "177: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 178: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 179: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B61").Value2 = "180: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 181: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 182: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B62").Value2 = "183: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 184: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 185: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B63").Value2 = "186: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 187: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 188: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B64").Value2 = "189: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 190: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 191: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B65").Value2 = "192: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 193: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 194: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B66").Value2 = "195: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 196: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 197: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B67").Value2 = "198: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 199: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 200: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B68").Value2 = "201: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 202: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 203: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B69").Value2 = "204: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 205: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 206: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B70").Value2 = "207: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 208: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 209: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B71").Value2 = "210: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 211: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 212: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B72").Value2 = "213: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 214: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 215: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B73").Value2 = "216: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 217: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 218: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B74").Value2 = "219: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 220: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 221: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B75").Value2 = "222: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 223: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 224: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B76").Value2 = "225: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 226: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 227: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B77").Value2 = "228: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 229: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 230: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B78").Value2 = "231: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 232: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 233: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B79").Value2 = "234: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 235: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 236: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B80").Value2 = "237: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 238: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 239: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B81").Value2 = "240: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 241: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 242: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B82").Value2 = "243: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 244: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 245: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B83").Value2 = "246: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 247: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 248: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B84").Value2 = "249: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 250: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 251: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B85").Value2 = "252: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 253: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 254: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B86").Value2 = "255: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 256: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 257: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B87").Value2 = "258: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 259: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 260: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B88").Value2 = "261: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 262: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 263: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B89").Value2 = "264: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 265: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 266: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B90").Value2 = "267: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 268: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 269: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B91").Value2 = "270: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 271: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 272: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B92").Value2 = "273: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 274: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 275: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B93").Value2 = "276: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 277: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 278: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B94").Value2 = "279: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 280: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 281: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B95").Value2 = "282: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 283: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 284: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B96").Value2 = "285: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 286: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 287: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B97").Value2 = "288: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 289: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 290: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B98").Value2 = "291: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 292: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 293: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B99").Value2 = "294: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 295: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 296: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B100").Value2 = "297: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 298: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 299: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B101").Value2 = "300: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 301: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 302: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B102").Value2 = "303: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 304: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 305: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B103").Value2 = "306: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 307: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 308: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B104").Value2 = "309: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 310: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 311: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B105").Value2 = "312: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 313: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 314: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B106").Value2 = "315: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 316: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 317: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B107").Value2 = "318: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 319: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 320: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B108").Value2 = "321: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 322: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 323: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B109").Value2 = "324: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 325: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 326: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B110").Value2 = "327: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 328: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 329: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B111").Value2 = "330: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 331: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 332: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B112").Value2 = "333: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 334: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 335: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B113").Value2 = "336: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 337: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 338: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B114").Value2 = "339: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 340: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 341: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B115").Value2 = "342: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 343: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 344: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B116").Value2 = "345: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 346: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 347: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B117").Value2 = "348: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 349: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 350: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B118").Value2 = "351: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 352: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 353: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B119").Value2 = "354: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 355: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 356: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
$ws.Range("B120").Value2 = "357: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 358: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] 359: [Sun Moon Asc Mercury Venus Mars Jupiter Saturn Uranus Neptune Pluto] "
